$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BRUNO's balance (row 3, column C) from 171035.75 to 40035.75
$ws.Cells.Item(3, 3).Value = 40035.75

# Delete the 15 rows for IRON..MONICA (originally rows 4-18)
$ws.Range("A4:A18").EntireRow.Delete()

# Delete the CAMILA row (originally row 20, now row 5 after the deletion above)
$ws.Range("A5").EntireRow.Delete()
